# Add 2022-Q3 data
#
# 1) Insert a brand-new worksheet named "2022-Q3" right after "总计" and
#    before "2022-Q2", populated with the per-fund holdings table for the
#    new quarter.
# 2) Update the "总计" (summary) sheet: insert a new top row for 2022-Q3
#    and push the existing quarterly rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: "总计" summary sheet — shift rows 2..7 down to 3..8, and write
# the new 2022-Q3 row into row 2.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Make room: duplicate formatting of the last existing data row (row 7)
# into the brand-new row 8, then fill rows top-to-bottom from the bottom
# up so we never clobber a value before it has been copied down.
$total.Range("A7:D7").Copy()
$total.Range("A8:D8").PasteSpecial(-4122)

# Row 8 <- old row 7 ("2021-Q1")
$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 5
$total.Range("D8").Value = 0.9

# Row 7 <- old row 6 ("2021-Q2")
$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 3
$total.Range("D7").Value = 0.88

# Row 6 <- old row 5 ("2021-Q3")
$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 2
$total.Range("D6").Value = 0.03

# Row 5 <- old row 4 ("2021-Q4")
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 1
$total.Range("D5").Value = 0.05

# Row 4 <- old row 3 ("2022-Q1")
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.08

# Row 3 <- old row 2 ("2022-Q2")
$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 0.95

# Row 2 <- new ("2022-Q3")
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 9
$total.Range("D2").Value = 0.37

# ---------------------------------------------------------------------
# Step 2: brand-new "2022-Q3" worksheet, inserted right after "总计".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($null, $total)
$newSheet.Name = "2022-Q3"

# Clone the header-row + data-row formatting from the "2022-Q2" sheet so
# style indices (bold/centered header, bordered index column) line up.
$q2.Range("A1:H1").Copy()
$newSheet.Range("A1:H1").PasteSpecial(-4122)
$q2.Range("A2:H2").Copy()
$newSheet.Range("A2:H10").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B:G hold text in the source data (fund codes with leading
# zeros, and numbers-as-text like "6.14"); force text format so Excel
# doesn't silently coerce them to numeric values.
$newSheet.Range("B2:G10").NumberFormat = "@"

function Set-FundRow($ws, $row, $idx, $code, $fundName, $scale, $stockPos, $posPct, $mv, $rank) {
    $ws.Range("A$row").Value = $idx
    $ws.Range("B$row").Value = $code
    $ws.Range("C$row").Value = $fundName
    $ws.Range("D$row").Value = $scale
    $ws.Range("E$row").Value = $stockPos
    $ws.Range("F$row").Value = $posPct
    $ws.Range("G$row").Value = $mv
    $ws.Range("H$row").Value = $rank
}

Set-FundRow $newSheet 2 0 "002666" "前海开源沪港深创新成长灵活配置混合A"   "6.14" "89.54" "3.59" "0.2204" 9
Set-FundRow $newSheet 3 1 "002667" "前海开源沪港深创新成长灵活配置混合C"   "2.27" "89.54" "3.59" "0.0815" 9
Set-FundRow $newSheet 4 2 "004099" "前海开源沪港深景气行业精选灵活配置混合" "0.46" "91.05" "4.79" "0.0220" 10
Set-FundRow $newSheet 5 3 "970021" "信达价值精选一年持有期灵活配置混合B"   "0.39" "51.28" "4.92" "0.0192" 5
Set-FundRow $newSheet 6 4 "009658" "汇丰晋信中小盘低波动策略股票A"         "0.85" "90.14" "1.92" "0.0163" 9
Set-FundRow $newSheet 7 5 "970020" "信达价值精选一年持有期灵活配置混合A"   "0.10" "51.28" "4.92" "0.0049" 5
Set-FundRow $newSheet 8 6 "004680" "前海开源裕瑞混合A"                     "0.37" "39.41" "1.08" "0.0040" 10
Set-FundRow $newSheet 9 7 "006190" "前海开源裕瑞混合C"                     "0.13" "39.41" "1.08" "0.0014" 10
Set-FundRow $newSheet 10 8 "009775" "汇丰晋信中小盘低波动策略股票C"        "0.04" "90.14" "1.92" "0.0008" 9
